# Apply the edit described by the commit:
# "added functionallaty to read in symbol table from excell"
#
# 1. Change the selection on the DB2 sheet (B7 -> C5).
# 2. Insert a new "Symbol_table" worksheet right before "READ_ME" and
#    fill it with a small sample symbol table (Symbol / Adress /
#    Data type / Comment).
# 3. Leave the selection/active-cell on the new sheet at D9, which makes
#    it the active tab (matching activeTab="4" in workbook.xml) and
#    clears the "tabSelected" flag that used to live on info_PLC /
#    READ_ME.

$wb = $excel.ActiveWorkbook

# --- 1. Update selection on DB2 ---
$db2 = $wb.Worksheets.Item("DB2")
$db2.Range("C5").Select()

# --- 2. Insert the new Symbol_table sheet before READ_ME ---
$readme = $wb.Worksheets.Item("READ_ME")
$symTable = $wb.Worksheets.Add($readme)
$symTable.Name = "Symbol_table"

# Fill the data rows first (this is the order that produced the shared
# string table ordering seen in the target workbook).
$symTable.Range("A2").Value = "test_name"
$symTable.Range("B2").Value = "I       0.0"
$symTable.Range("C2").Value = "BOOL"

$symTable.Range("A3").Value = "other name"
$symTable.Range("B3").Value = "I       0.1"
$symTable.Range("C3").Value = "BOOL"

# Then the header row.
$symTable.Range("A1").Value = "Symbol"
$symTable.Range("B1").Value = "Adress"
$symTable.Range("C1").Value = "Data type"
$symTable.Range("D1").Value = "Comment"

# Then one more data row, added afterwards.
$symTable.Range("A4").Value = "out"
$symTable.Range("B4").Value = "Q       0.1"
$symTable.Range("C4").Value = "BOOL"
$symTable.Range("D4").Value = "test comment"

# --- 3. Leave the selection / active sheet as in the target file ---
$symTable.Range("D9").Select()
